$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.474.10'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '1.813.24'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('E6').Value = '  +2.99%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '38.37'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.06%  '
$ws.Range('E9').Value = '  -3.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0676'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0973'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').Value = '2.074.77'
$ws.Range('E12').Value = '  +0.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.23'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.19%  '
$ws.Range('D14').Value = '1.829.20'
$ws.Range('E14').Value = '  +1.41%  '
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').Value = '34.462.20'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.33'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.28'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').Value = '0.0₃0775'
$ws.Range('E20').Value = '  -2.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('E23').Value = '  -1.35%  '
$ws.Range('E24').Value = '  +3.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '170.34'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.82'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.58'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.31%  '
$ws.Range('E28').Value = '  +2.36%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.80'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.06%  '
$ws.Range('E32').Value = '  -2.52%  '
$ws.Range('E33').Value = '  -4.81%  '
$ws.Range('E34').Value = '  +0.24%  '
$ws.Range('D35').Value = '1.360.52'
$ws.Range('E35').Value = '  -2.41%  '
$ws.Range('E36').Value = '  -4.00%  '
$ws.Range('E37').Value = '  -0.80%  '
$ws.Range('E38').Value = '  -4.66%  '
$ws.Range('E39').Value = '  -1.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.44'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.28%  '
$ws.Range('E41').Value = '  -1.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.953'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.88%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '81.85'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.38%  '
$ws.Range('E44').Value = '  -0.71%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.80'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.07%  '
$ws.Range('E46').Value = '  +1.35%  '
$ws.Range('D47').Value = '1.975.73'
$ws.Range('E47').Value = '  +0.56%  '
$ws.Range('E48').Value = '  -4.59%  '
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '102.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.58%  '
$ws.Range('E51').Value = '  -5.06%  '
